# Generate Report for Handoff
# Adds a new file entry (a29441b6-aa6f-4c64-8ced-f0e96db3b680) as row 3
# on the Overview, zh-cn and de-de sheets, mirroring the existing
# 87295fab-7e21-42f7-81c6-2353ba2e6415 row.

$wb = $excel.ActiveWorkbook

$newGuid = "a29441b6-aa6f-4c64-8ced-f0e96db3b680"
$zhHash = "70de58809b20a8d8ab75317e2c0a9e08ebe0d72e"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/10d8aeef3c22dcce1921eec011f99d652ca29666/e2e/$newGuid.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93612d475839679fe25979baad56ac776bff1997/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$zhHash.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad6511544582fa5a3052925d35a117916e269c8b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$zhHash.de-de.xlf"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "$newGuid.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-31-21 00:31:56"

$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, "", "", "$newGuid.md")

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "$newGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "$newGuid.$zhHash.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 00:31:53"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $zhXlfUrl, "", "", "$newGuid.$zhHash.zh-cn.xlf")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "$newGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "$newGuid.$zhHash.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 00:31:56"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $deXlfUrl, "", "", "$newGuid.$zhHash.de-de.xlf")
